$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts old F -> G)
$ws.Columns("F:F").Insert()

# New header for the inserted column F
$ws.Range("G1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Index (ČR=100)"

# New data values for column F (rows 2-15)
$values = @(249, 95, 92, 91, 73, 73, 72, 72, 72, 71, 71, 71, 70, 69)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
